$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "平潭发展"
$ws.Range("B2").Value = "平潭发展"
$ws.Range("C2").Value = "中际旭创"
$ws.Range("A3").Value = "工业富联"
$ws.Range("B3").Value = "山子高科"
$ws.Range("C3").Value = "平潭发展"
$ws.Range("A4").Value = "神州信息"
$ws.Range("B4").Value = "隆基绿能"
$ws.Range("C4").Value = "福龙马"
$ws.Range("A5").Value = "多氟多"
$ws.Range("B5").Value = "工业富联"
$ws.Range("C5").Value = "和而泰"
$ws.Range("A6").Value = "隆基绿能"
$ws.Range("B6").Value = "天齐锂业"
$ws.Range("C6").Value = "三花智控"
$ws.Range("A7").Value = "天际股份"
$ws.Range("B7").Value = "神州信息"
$ws.Range("C7").Value = "山子高科"
$ws.Range("A8").Value = "山子高科"
$ws.Range("B8").Value = "鹏辉能源"
$ws.Range("C8").Value = "神州信息"
$ws.Range("A9").Value = "天齐锂业"
$ws.Range("B9").Value = "赣锋锂业"
$ws.Range("C9").Value = "赛力斯"
$ws.Range("A10").Value = "格尔软件"
$ws.Range("B10").Value = "多氟多"
$ws.Range("C10").Value = "格尔软件"
$ws.Range("A11").Value = "福龙马"
$ws.Range("B11").Value = "北方稀土"
$ws.Range("C11").Value = "隆基绿能"
$ws.Range("A12").Value = "北方稀土"
$ws.Range("B12").Value = "阳光电源"
$ws.Range("C12").Value = "江特电机"
$ws.Range("A13").Value = "江特电机"
$ws.Range("B13").Value = "天际股份"
$ws.Range("C13").Value = "士兰微"
$ws.Range("A14").Value = "赣锋锂业"
$ws.Range("B14").Value = "江特电机"
$ws.Range("C14").Value = "天际股份"
$ws.Range("A15").Value = "鹏辉能源"
$ws.Range("B15").Value = "江波龙"
$ws.Range("C15").Value = "常铝股份"
$ws.Range("A16").Value = "三花智控"
$ws.Range("B16").Value = "格尔软件"
$ws.Range("C16").Value = "天齐锂业"
$ws.Range("A17").Value = "科大国创"
$ws.Range("B17").Value = "三花智控"
$ws.Range("C17").Value = "工业富联"
$ws.Range("A18").Value = "国盾量子"
$ws.Range("B18").Value = "东方财富"
$ws.Range("C18").Value = "北方稀土"
$ws.Range("A19").Value = "阳光电源"
$ws.Range("B19").Value = "福龙马"
$ws.Range("C19").Value = "多氟多"
$ws.Range("A20").Value = "新易盛"
$ws.Range("B20").Value = "科大国创"
$ws.Range("C20").Value = "黄河旋风"
$ws.Range("A21").Value = "和而泰"
$ws.Range("B21").Value = "新易盛"
$ws.Range("C21").Value = "赣锋锂业"
